# Apply crypto price/volume updates, including two row pair swaps (16/17, 22/23, 41/42)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-CellText "D2" "29.383.10"
Set-CellText "E2" "  -0.42%  "

# Row 3
Set-CellText "D3" "1.847.57"
Set-CellText "E3" "  -0.13%  "

# Row 4
Set-CellText "E4" "  -0.12%  "

# Row 5
Set-CellText "D5" "241.23"
Set-CellText "E5" "  -0.85%  "

# Row 6
Set-CellText "D6" "0.6257"
Set-CellText "E6" "  -3.85%  "

# Row 7
Set-CellText "D7" "1.001"
Set-CellText "E7" "  -0.05%  "

# Row 8
Set-CellText "D8" "0.07614"
Set-CellText "E8" "  +1.57%  "

# Row 9
Set-CellText "D9" "0.2969"
Set-CellText "E9" "  -0.24%  "

# Row 10
Set-CellText "D10" "24.41"
Set-CellText "E10" "  -0.07%  "

# Row 11
Set-CellText "D11" "2.057.45"
Set-CellText "E11" "  +11.03%  "

# Row 12
Set-CellText "D12" "0.07721"
Set-CellText "E12" "  +1.06%  "

# Row 13
Set-CellText "D13" "4.989"
Set-CellText "E13" "  -0.88%  "

# Row 14
Set-CellText "D14" "0.6875"
Set-CellText "E14" "  +0.29%  "

# Row 15
Set-CellText "D15" "82.92"
Set-CellText "E15" "  -0.75%  "

# Row 16
Set-CellText "B16" "WrappedliquidstakedEther2.0"
Set-CellText "C16" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText "D16" "2.286.12"
Set-CellText "E16" "  +8.33%  "

# Row 17
Set-CellText "B17" "ShibaInu"
Set-CellText "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText "D17" "0.000009943"
Set-CellText "E17" "  +4.26%  "

# Row 18
Set-CellText "D18" "6.152"
Set-CellText "E18" "  +0.37%  "

# Row 19
Set-CellText "D19" "29.632.31"

# Row 20
Set-CellText "D20" "230.73"
Set-CellText "E20" "  -2.69%  "

# Row 21
Set-CellText "E21" "  -0.63%  "

# Row 22
Set-CellText "B22" "Chainlink"
Set-CellText "C22" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellText "D22" "7.717"
Set-CellText "E22" "  +0.28%  "

# Row 23
Set-CellText "B23" "Dai"
Set-CellText "C23" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText "D23" "1.000"
Set-CellText "E23" "  -0.08%  "

# Row 24
Set-CellText "E24" "  -0.12%  "

# Row 25
Set-CellText "D25" "154.65"
Set-CellText "E25" "  -1.76%  "

# Row 26
Set-CellText "D26" "0.1388"
Set-CellText "E26" "  -2.19%  "

# Row 27
Set-CellText "D27" "8.463"
Set-CellText "E27" "  -0.69%  "

# Row 28
Set-CellText "E28" "  -0.88%  "

# Row 29
Set-CellText "E29" "  -1.01%  "

# Row 30
Set-CellText "D30" "0.05813"
Set-CellText "E30" "  -4.32%  "

# Row 31
Set-CellText "D31" "1.256"
Set-CellText "E31" "  -0.12%  "

# Row 32
Set-CellText "D32" "4.122"

# Row 33
Set-CellText "E33" "  -1.40%  "

# Row 34
Set-CellText "D34" "1.865"
Set-CellText "E34" "  +0.01%  "

# Row 35
Set-CellText "E35" "  -2.30%  "

# Row 36
Set-CellText "D36" "0.7172"
Set-CellText "E36" "  -1.22%  "

# Row 37
Set-CellText "D37" "2.600"
Set-CellText "E37" "  -0.06%  "

# Row 38
Set-CellText "D38" "1.250.45"
Set-CellText "E38" "  +4.05%  "

# Row 39
Set-CellText "D39" "2.792"
Set-CellText "E39" "  -0.47%  "

# Row 40
Set-CellText "E40" "  +0.98%  "

# Row 41
Set-CellText "B41" "TrustWalletToken"
Set-CellText "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText "D41" "0.9090"
Set-CellText "E41" "  -0.09%  "

# Row 42
Set-CellText "B42" "RocketPoolETH"
Set-CellText "C42" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-CellText "D42" "2.193.74"
Set-CellText "E42" "  +8.77%  "

# Row 43
Set-CellText "D43" "6.072"
Set-CellText "E43" "  -2.42%  "

# Row 44
Set-CellText "D44" "0.9997"
Set-CellText "E44" "  -0.11%  "

# Row 45
Set-CellText "D45" "101.92"
Set-CellText "E45" "  +0.30%  "

# Row 46
Set-CellText "D46" "67.45"
Set-CellText "E46" "  +1.39%  "

# Row 47
Set-CellText "D47" "7.327"
Set-CellText "E47" "  -0.80%  "

# Row 48
Set-CellText "E48" "  -2.15%  "

# Row 49
Set-CellText "D49" "9.144"
Set-CellText "E49" "  -0.05%  "

# Row 50
Set-CellText "D50" "0.4024"
Set-CellText "E50" "  -0.63%  "

# Row 51
Set-CellText "D51" "1.699"
Set-CellText "E51" "  +2.56%  "
